# Daily attendance processing - 2025-12-09 19:24:07
# Normalizes the "Recorded By" column (G) so that when a cell holds more
# than one comma-separated name and "System" isn't already first, the
# first two names are swapped (bringing "System" to the front where
# present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0] -ne "System") {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
